# Appends newly received loss-of-sale records (rows 41-51 / records #39-#49)
# to the Walk-In Report sheet, keeping the existing schema (cols A-K).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRecords = @(
    @(41, 39, "'22-12-2025", "SHIBIL", 8891013925, "'31-12-2025", "MUHAMMED YASHEER M A", "Loss", "CUSTOMER INTERNAL ISSUES", "FAMILY DISAPPROVEL", "-", "nokeet paraya parannu"),
    @(42, 40, "'22-12-2025", "sharath", 9526271630, "'15-02-2026", "ASHIK A", "Loss", "ENQUIRY", "ENQUIRY WITHOUT BRIDE/FAMILY", "-", "next week varam"),
    @(43, 41, "'23-12-2025", "BHARATH", 9020409040, "'29-12-2025", "MUHAMMED YASHEER M A", "Loss", "SIZE NOT SUITABLE", "SIZE TOO LARGE", "-", "size 44 venam athil ee shopil collections illathath kond eduthilla"),
    @(44, 42, "'23-12-2025", "Prince", 7034191408, "'28-01-2026", "NAVAS A", "Loss", "CUSTOMER INTERNAL ISSUES", "FAMILY DISAPPROVEL", "-", "family chothichite varam parannu"),
    @(45, 43, "'24-12-2025", "krishadas", 859006813, "'18-01-2026", "NAVAS A", "Loss", "ENQUIRY", "ENQUIRY WITHOUT TRIAL", "-", "groom not comming"),
    @(46, 44, "'25-12-2025", "vibin", 897343577, "'08-02-2026", "NAVAS A", "Loss", "ENQUIRY", "ENQUIRY WITHOUT TRIAL", "-", "bride dress eduthitila edutite varam parannu"),
    @(47, 45, "'25-12-2025", "MANIKANDHAN", 8075294343, "'24-01-2026", "MUHAMMED YASHEER M A", "Loss", "CUSTOMER INTERNAL ISSUES", "FAMILY DISAPPROVEL", "-", "nokeet paraya parannu"),
    @(48, 46, "'25-12-2025", "RIYAS", 9539111756, "'18-01-2026", "MUHAMMED YASHEER M A", "Loss", "CUSTOMER INTERNAL ISSUES", "FAMILY DISAPPROVEL", "-", "nokeet paraya parannu"),
    @(49, 47, "'25-12-2025", "RAFEEQ", 8848314475, "'04-01-2026", "MUHAMMED YASHEER M A", "Loss", "CUSTOMER INTERNAL ISSUES", "FAMILY DISAPPROVEL", "-", "next week ayitt vara parannu"),
    @(50, 48, "'25-12-2025", "SURESH", 9074734689, "'25-01-2026", "MUHAMMED YASHEER M A", "Loss", "CUSTOMER INTERNAL ISSUES", "FAMILY DISAPPROVEL", "-", "nokeet paraya parannu"),
    @(51, 49, "'25-12-2025", "AMEESH", 9809016325, "'01-02-2026", "NAVAS A", "Loss", "CUSTOMER INTERNAL ISSUES", "FAMILY DISAPPROVEL", "-", "long date ann family ayitt choyichitt paraya parannu"),
)

foreach ($rec in $newRecords) {
    $r  = $rec[0]
    $ws.Cells.Item($r, 1).Value  = $rec[1]   # A: # (serial number)
    $ws.Cells.Item($r, 2).Value  = $rec[2]   # B: Date (text)
    $ws.Cells.Item($r, 3).Value  = $rec[3]   # C: Customer Name
    $ws.Cells.Item($r, 4).Value  = $rec[4]   # D: Contact
    $ws.Cells.Item($r, 5).Value  = $rec[5]   # E: Function Date (text)
    $ws.Cells.Item($r, 6).Value  = $rec[6]   # F: Staff
    $ws.Cells.Item($r, 7).Value  = $rec[7]   # G: Status
    $ws.Cells.Item($r, 8).Value  = $rec[8]   # H: Category
    $ws.Cells.Item($r, 9).Value  = $rec[9]   # I: Sub Category
    $ws.Cells.Item($r, 10).Value = $rec[10]  # J: Repeat count
    $ws.Cells.Item($r, 11).Value = $rec[11]  # K: Remarks

    # Columns A and D are numeric in the existing data, formatted as plain integers.
    $ws.Cells.Item($r, 1).NumberFormat = "0"
    $ws.Cells.Item($r, 4).NumberFormat = "0"
}
